# Updated cryptos list on Fri Oct 20 10:51:07 UTC 2023 with GitHub Actions
# Refresh the coin ranking table with the latest scraped Price / Volume(1h)
# values (and re-order two pairs of rows whose relative ranking flipped).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells D5,D6,D8,D10,D11,D15-D30(partial),D32,D36,D37,D40,D41,D43,D45-D47,D49
# hold plain-decimal price strings (e.g. "214.55"). Excel's COM Value setter
# auto-coerces such strings to numbers, which would lose the original text
# formatting (e.g. trailing zeros like "0.0600" -> 0.06) and introduce
# floating point noise. Forcing NumberFormat to Text ("@") first keeps them
# stored as text, matching the source data feed.

$ws.Range('D2').Value = '29.814.23'
$ws.Range('E2').Value = '  +4.85%  '
$ws.Range('D3').Value = '1.614.18'
$ws.Range('E3').Value = '  +4.13%  '
$ws.Range('E4').Value = '  -0.41%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.55'
$ws.Range('E5').Value = '  +1.88%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.517'
$ws.Range('E6').Value = '  +7.07%  '
$ws.Range('E7').Value = '  -0.50%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '26.94'
$ws.Range('E8').Value = '  +12.55%  '
$ws.Range('E9').Value = '  +3.76%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0600'
$ws.Range('E10').Value = '  +2.92%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0914'
$ws.Range('E11').Value = '  +2.91%  '
$ws.Range('D12').Value = '1.845.09'
$ws.Range('E12').Value = '  +4.15%  '
$ws.Range('D13').Value = '1.614.38'
$ws.Range('E13').Value = '  +4.25%  '
$ws.Range('B14').Value = 'WrappedBTC'
$ws.Range('C14').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D14').Value = '29.851.63'
$ws.Range('E14').Value = '  +5.02%  '
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.542'
$ws.Range('E15').Value = '  +6.52%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.76'
$ws.Range('E16').Value = '  +3.92%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '246.02'
$ws.Range('E17').Value = '  +7.76%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '63.66'
$ws.Range('E18').Value = '  +4.44%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.63'
$ws.Range('E19').Value = '  +4.37%  '
$ws.Range('D20').Value = '0.0₃0697'
$ws.Range('E20').Value = '  +3.69%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.994'
$ws.Range('E21').Value = '  -0.48%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.06'
$ws.Range('E22').Value = '  +4.52%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.29'
$ws.Range('E23').Value = '  +4.52%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.11'
$ws.Range('E24').Value = '  +4.66%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '155.81'
$ws.Range('E25').Value = '  +3.05%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.41'
$ws.Range('E26').Value = '  +4.58%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.108'
$ws.Range('E27').Value = '  +5.60%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.42'
$ws.Range('E28').Value = '  +3.49%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.996'
$ws.Range('E29').Value = '  -0.38%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0473'
$ws.Range('E30').Value = '  +1.65%  '
$ws.Range('E31').Value = '  +0.78%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.25'
$ws.Range('E32').Value = '  +3.22%  '
$ws.Range('D33').Value = '1.442.60'
$ws.Range('E33').Value = '  +4.29%  '
$ws.Range('E34').Value = '  +4.28%  '
$ws.Range('E35').Value = '  -1.44%  '
$ws.Range('B36').Value = 'MXToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.85'
$ws.Range('E36').Value = '  +10.97%  '
$ws.Range('B37').Value = 'LidoDAOToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.52'
$ws.Range('E37').Value = '  +3.24%  '
$ws.Range('E38').Value = '  +0.48%  '
$ws.Range('E39').Value = '  +3.32%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '56.49'
$ws.Range('E40').Value = '  +31.05%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.540'
$ws.Range('E41').Value = '  +5.89%  '
$ws.Range('E42').Value = '  +2.20%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.800'
$ws.Range('E43').Value = '  +3.89%  '
$ws.Range('E44').Value = '  -0.46%  '
$ws.Range('B45').Value = 'Kaspa'
$ws.Range('C45').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0470'
$ws.Range('E45').Value = '  +2.64%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '66.60'
$ws.Range('E46').Value = '  +7.87%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '5.35'
$ws.Range('E47').Value = '  +0.53%  '
$ws.Range('D48').Value = '1.754.86'
$ws.Range('E48').Value = '  +4.31%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '86.99'
$ws.Range('E49').Value = '  +2.13%  '
$ws.Range('E50').Value = '  -4.42%  '
$ws.Range('D51').Value = '0.0₆0102'
$ws.Range('E51').Value = '  +1.13%  '
